$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1610.7858
$ws.Range("J17").Value = 1260.7037
$ws.Range("L17").Value = 3782.1111
$ws.Range("N17").Value = -4118.1111
$ws.Range("H113").Value = 3659.5
$ws.Range("I113").Value = 3849.375
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 3849.375
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -595.375
$ws.Range("N113").Value = -9408
$ws.Range("H132").Value = 1327.8889
$ws.Range("I132").Value = 1393.6
$ws.Range("K132").Value = 4180.799999999999
$ws.Range("M132").Value = -1650.799999999999
$ws.Range("H137").Value = 2388.8518
$ws.Range("I137").Value = 2145.8333
$ws.Range("K137").Value = 6437.499899999999
$ws.Range("M137").Value = -3887.499899999999
$ws.Range("H140").Value = 77497.336
$ws.Range("J140").Value = 77497.336
$ws.Range("L140").Value = 77497.336
$ws.Range("N140").Value = -87857.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 3668
$ws.Range("I35").Value = 3668
$ws.Range("K35").Value = 3668
$ws.Range("M35").Value = -3262
$ws.Range("H132").Value = 1884.7916
$ws.Range("I132").Value = 1726
$ws.Range("K132").Value = 5178
$ws.Range("M132").Value = -2648

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 199.5
$ws.Range("I22").Value = 199.5
$ws.Range("K22").Value = 199.5
$ws.Range("M22").Value = -26.5
$ws.Range("H37").Value = 50000
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6250831.5
$ws.Range("J22").Value = 10417534
$ws.Range("L22").Value = 10417534
$ws.Range("N22").Value = -10418234
$ws.Range("H58").Value = 3345353.5
$ws.Range("I58").Value = 3953482.2
$ws.Range("K58").Value = 3953482.2
$ws.Range("M58").Value = -3953279.2
$ws.Range("H99").Value = 668414.4
$ws.Range("I99").Value = 1430044.1
$ws.Range("J99").Value = 1988.375
$ws.Range("K99").Value = 1430044.1
$ws.Range("L99").Value = 1988.375
$ws.Range("M99").Value = -1428546.1
$ws.Range("N99").Value = -4984.375
$ws.Range("H126").Value = 668414.4
$ws.Range("I126").Value = 1430044.1
$ws.Range("J126").Value = 1988.375
$ws.Range("K126").Value = 4290132.300000001
$ws.Range("L126").Value = 5965.125
$ws.Range("M126").Value = -4287662.300000001
$ws.Range("N126").Value = -10905.125
$ws.Range("H136").Value = 3345353.5
$ws.Range("I136").Value = 3953482.2
$ws.Range("K136").Value = 11860446.6
$ws.Range("M136").Value = -11857896.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2256.9795
$ws.Range("I68").Value = 650.6667
$ws.Range("J68").Value = 2361.739
$ws.Range("K68").Value = 1952.0001
$ws.Range("L68").Value = 7085.217000000001
$ws.Range("M68").Value = -1141.0001
$ws.Range("N68").Value = -8707.217000000001
$ws.Range("H71").Value = 2256.9795
$ws.Range("I71").Value = 650.6667
$ws.Range("J71").Value = 2361.739
$ws.Range("K71").Value = 5856.0003
$ws.Range("L71").Value = 21255.651
$ws.Range("M71").Value = -1800.0003
$ws.Range("N71").Value = -29367.651
$ws.Range("H107").Value = 1610.8334
$ws.Range("I107").Value = 1082.2222
$ws.Range("K107").Value = 3246.6666
$ws.Range("M107").Value = -1326.6666
$ws.Range("H115").Value = 4617.5454
$ws.Range("I115").Value = 1666.3334
$ws.Range("J115").Value = 5724.25
$ws.Range("K115").Value = 4999.0002
$ws.Range("L115").Value = 17172.75
$ws.Range("M115").Value = -3824.0002
$ws.Range("N115").Value = -19522.75
$ws.Range("H129").Value = 32273.088
$ws.Range("J129").Value = 52597.285
$ws.Range("L129").Value = 157791.855
$ws.Range("N129").Value = -167791.855
$ws.Range("H131").Value = 16153123
$ws.Range("J131").Value = 26601.715
$ws.Range("L131").Value = 79805.145
$ws.Range("N131").Value = -89885.145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 1674888.9
$ws.Range("I20").Value = 7500000
$ws.Range("J20").Value = 10571.429
$ws.Range("K20").Value = 7500000
$ws.Range("L20").Value = 10571.429
$ws.Range("M20").Value = -7499755
$ws.Range("N20").Value = -11061.429
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("I126").Value = 3475179
$ws.Range("J126").Value = 2385.3333
$ws.Range("K126").Value = 10425537
$ws.Range("L126").Value = 7155.999899999999
$ws.Range("M126").Value = -10423067
$ws.Range("N126").Value = -12095.9999
$ws.Range("H132").Value = 1926267.9
$ws.Range("I132").Value = 2406585
$ws.Range("K132").Value = 7219755
$ws.Range("M132").Value = -7217225

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 2817.3333
$ws.Range("J7").Value = 4732.3335
$ws.Range("K7").Value = 2817.3333
$ws.Range("L7").Value = 4732.3335
$ws.Range("M7").Value = -2705.3333
$ws.Range("N7").Value = -4956.3335
$ws.Range("H40").Value = 5163.5884
$ws.Range("I40").Value = 4520.5557
$ws.Range("J40").Value = 5887
$ws.Range("K40").Value = 4520.5557
$ws.Range("L40").Value = 5887
$ws.Range("M40").Value = -4384.5557
$ws.Range("N40").Value = -6159
$ws.Range("H63").Value = 12000
$ws.Range("J63").Value = 12000
$ws.Range("L63").Value = 12000
$ws.Range("N63").Value = -13498
$ws.Range("H66").Value = 12000
$ws.Range("J66").Value = 12000
$ws.Range("L66").Value = 36000
$ws.Range("N66").Value = -43488
$ws.Range("H122").Value = 5672.65
$ws.Range("I122").Value = 5050.4443
$ws.Range("J122").Value = 6181.727
$ws.Range("K122").Value = 15151.3329
$ws.Range("L122").Value = 18545.181
$ws.Range("M122").Value = -12701.3329
$ws.Range("N122").Value = -23445.181
$ws.Range("I126").Value = 2817.3333
$ws.Range("J126").Value = 4732.3335
$ws.Range("K126").Value = 8451.999899999999
$ws.Range("L126").Value = 14197.0005
$ws.Range("M126").Value = -5981.999899999999
$ws.Range("N126").Value = -19137.0005
$ws.Range("H136").Value = 4642.357
$ws.Range("I136").Value = 2689.2222
$ws.Range("J136").Value = 8158
$ws.Range("K136").Value = 8067.6666
$ws.Range("L136").Value = 24474
$ws.Range("M136").Value = -5517.6666
$ws.Range("N136").Value = -29574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 25866.666
$ws.Range("J70").Value = 25866.666
$ws.Range("L70").Value = 25866.666
$ws.Range("N70").Value = -26496.666
$ws.Range("H73").Value = 25866.666
$ws.Range("J73").Value = 25866.666
$ws.Range("L73").Value = 25866.666
$ws.Range("N73").Value = -28050.666
$ws.Range("H126").Value = 4443.857
$ws.Range("I126").Value = 3634.4666
$ws.Range("K126").Value = 10903.3998
$ws.Range("M126").Value = -8433.399800000001
